$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Price + Volume(1h) columns) to match the latest
# coinranking.com snapshot. D-column values are prefixed with a leading
# apostrophe so Excel stores them as literal text (preserving formats such
# as trailing zeros / "." thousands separators) instead of auto-converting
# them to numbers.

$ws.Range("D2").Value = "'40.669.58"
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("D3").Value = "'2.371.81"
$ws.Range("E3").Value = "  -4.14%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'310.77"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").Value = "'86.42"
$ws.Range("E6").Value = "  -6.47%  "
$ws.Range("E7").Value = "  -3.99%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("D11").Value = "'30.32"
$ws.Range("E11").Value = "  -8.66%  "
$ws.Range("D12").Value = "'0.109"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").Value = "'2.736.52"
$ws.Range("E13").Value = "  -4.16%  "
$ws.Range("D14").Value = "'6.53"
$ws.Range("E14").Value = "  -5.45%  "
$ws.Range("D15").Value = "'14.94"
$ws.Range("E15").Value = "  -3.56%  "
$ws.Range("D16").Value = "'2.396.75"
$ws.Range("E16").Value = "  -3.25%  "
$ws.Range("E17").Value = "  -4.35%  "
$ws.Range("D18").Value = "'40.613.46"
$ws.Range("E18").Value = "  -2.35%  "
$ws.Range("D19").Value = "'0.0₃0911"
$ws.Range("E19").Value = "  -3.61%  "
$ws.Range("E20").Value = "  -5.09%  "
$ws.Range("D21").Value = "'68.35"
$ws.Range("E21").Value = "  -3.36%  "
$ws.Range("E22").Value = "  -4.74%  "
$ws.Range("D23").Value = "'234.96"
$ws.Range("E23").Value = "  -2.39%  "
$ws.Range("D24").Value = "'2.58"
$ws.Range("E24").Value = "  -6.29%  "
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").Value = "'1.80"
$ws.Range("E26").Value = "  -7.50%  "
$ws.Range("D27").Value = "'23.78"
$ws.Range("E27").Value = "  -4.23%  "
$ws.Range("D28").Value = "'2.15"
$ws.Range("E28").Value = "  -3.37%  "
$ws.Range("D29").Value = "'9.23"
$ws.Range("E29").Value = "  -4.86%  "
$ws.Range("D30").Value = "'34.21"
$ws.Range("E30").Value = "  -6.84%  "
$ws.Range("D31").Value = "'153.65"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  -5.38%  "
$ws.Range("D34").Value = "'0.0727"
$ws.Range("E34").Value = "  -4.99%  "
$ws.Range("E35").Value = "  -5.86%  "
$ws.Range("E36").Value = "  -2.44%  "
$ws.Range("D37").Value = "'15.99"
$ws.Range("E37").Value = "  -7.53%  "
$ws.Range("E38").Value = "  -4.05%  "
$ws.Range("E39").Value = "  -4.55%  "
$ws.Range("E40").Value = "  -8.33%  "
$ws.Range("D41").Value = "'3.83"
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("D42").Value = "'2.40"
$ws.Range("E42").Value = "  -3.44%  "
$ws.Range("D43").Value = "'1.956.62"
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("E44").Value = "  -4.96%  "
$ws.Range("D45").Value = "'17.63"
$ws.Range("E45").Value = "  -6.24%  "
$ws.Range("D46").Value = "'9.36"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("E47").Value = "  -9.97%  "
$ws.Range("D48").Value = "'2.599.74"
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("D49").Value = "'92.97"
$ws.Range("E49").Value = "  -5.16%  "
$ws.Range("D50").Value = "'71.70"
$ws.Range("E50").Value = "  -5.33%  "
$ws.Range("D51").Value = "'50.12"
$ws.Range("E51").Value = "  -4.32%  "
